$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.483.94"
$ws.Range("E2").Value = "  +3.10%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.605.96"
$ws.Range("E3").Value = "  +2.73%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'212.65"
$ws.Range("E5").Value = "  +1.03%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.521"
$ws.Range("E6").Value = "  +6.72%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'26.99"
$ws.Range("E8").Value = "  +7.54%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'43.55"
$ws.Range("E9").Value = "  -0.70%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  +2.58%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.34%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.0911"
$ws.Range("E12").Value = "  +1.61%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "1.834.88"
$ws.Range("E13").Value = "  +2.69%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.622.23"
$ws.Range("E14").Value = "  +3.70%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "29.482.70"
$ws.Range("E15").Value = "  +3.07%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "'0.536"
$ws.Range("E16").Value = "  +4.11%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  +1.89%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'63.16"
$ws.Range("E18").Value = "  +3.05%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'241.41"
$ws.Range("E19").Value = "  +5.22%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'7.63"
$ws.Range("E20").Value = "  +3.81%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +1.97%  "

# Row 22 - Dai
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.04%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +2.62%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  +2.42%  "

# Row 25
$ws.Range("E25").Value = "  +0.59%  "

# Row 26
$ws.Range("D26").Value = "'154.79"
$ws.Range("E26").Value = "  +2.52%  "

# Row 27
$ws.Range("E27").Value = "  +5.06%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'15.28"
$ws.Range("E28").Value = "  +3.48%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "'6.39"
$ws.Range("E29").Value = "  +2.55%  "

# Row 30 - BinanceUSD
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.00%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +2.67%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +1.12%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'3.23"
$ws.Range("E33").Value = "  +1.89%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +4.30%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.414.77"
$ws.Range("E35").Value = "  +1.99%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +0.65%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  +3.40%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  +4.92%  "

# Row 39 - HuobiToken
$ws.Range("E39").Value = "  +0.41%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +2.56%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  +3.53%  "

# Row 42 - RenderToken
$ws.Range("E42").Value = "  +0.95%  "

# Row 43 - Kaspa
$ws.Range("D43").Value = "'0.0488"
$ws.Range("E43").Value = "  +6.17%  "

# Row 44 - ARBITRUM
$ws.Range("E44").Value = "  +3.25%  "

# Row 45 - was BitcoinSV, now PaxDollar
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.03%  "

# Row 46 - was PaxDollar, now BitcoinSV
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "'52.89"
$ws.Range("E46").Value = "  +22.47%  "

# Row 47 - Aave
$ws.Range("D47").Value = "'65.67"
$ws.Range("E47").Value = "  +2.75%  "

# Row 48 - FraxShare
$ws.Range("E48").Value = "  +1.09%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "1.746.78"
$ws.Range("E49").Value = "  +2.94%  "

# Row 50 - WEMIXToken
$ws.Range("D50").Value = "'0.857"
$ws.Range("E50").Value = "  -1.36%  "

# Row 51 - Quant
$ws.Range("D51").Value = "'86.84"
$ws.Range("E51").Value = "  +1.96%  "
